# Update countries & provincias Spain
# Applies the 27-Jun-2020 02:07 COVID-19 data refresh:
#  - several countries swap rank (same row, new country name + new stats)
#  - numeric case/death counters refreshed for ~40 rows
#  - footer timestamp updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (A1)
$ws.Range('A1').Value = 'Datos actualizados a 27 de Junio de 2020 a las 02:07'

# Row 4: Estados Unidos
$ws.Range('B4').Value = 2552446
$ws.Range('C4').Value = 46831
$ws.Range('D4').Value = 1065516
$ws.Range('E4').Value = 1359311
$ws.Range('G4').Value = 642
$ws.Range('H4').Value = 127619

# Row 5: Brasil
$ws.Range('B5').Value = 1280054
$ws.Range('C5').Value = 46907
$ws.Range('E5').Value = 526419
$ws.Range('G5').Value = 1055
$ws.Range('H5').Value = 56109

# Row 22: Canada
$ws.Range('B22').Value = 102794
$ws.Range('C22').Value = 172
$ws.Range('D22').Value = 65726
$ws.Range('E22').Value = 28560
$ws.Range('G22').Value = 4
$ws.Range('H22').Value = 8508

# Row 30: Argentina
$ws.Range('A30').Value = 'Argentina'
$ws.Range('B30').Value = 55343
$ws.Range('C30').Value = 2886
$ws.Range('D30').Value = 18416
$ws.Range('E30').Value = 35743
$ws.Range('G30').Value = 34
$ws.Range('H30').Value = 1184

# Row 31: Ecuador
$ws.Range('A31').Value = 'Ecuador'
$ws.Range('B31').Value = 53856
$ws.Range('C31').Value = 700
$ws.Range('D31').Value = 26493
$ws.Range('E31').Value = 22957
$ws.Range('G31').Value = 63
$ws.Range('H31').Value = 4406

# Row 51: Nigeria
$ws.Range('A51').Value = 'Nigeria'
$ws.Range('B51').Value = 23298
$ws.Range('C51').Value = 684
$ws.Range('D51').Value = 8253
$ws.Range('E51').Value = 14491
$ws.Range('G51').Value = 5
$ws.Range('H51').Value = 554

# Row 52: Armenia
$ws.Range('A52').Value = 'Armenia'
$ws.Range('B52').Value = 23247
$ws.Range('C52').Value = 759
$ws.Range('D52').Value = 12149
$ws.Range('E52').Value = 10688
$ws.Range('G52').Value = 13
$ws.Range('H52').Value = 410

# Row 53: Israel
$ws.Range('A53').Value = 'Israel'
$ws.Range('B53').Value = 22800
$ws.Range('C53').Value = 400
$ws.Range('D53').Value = 16872
$ws.Range('E53').Value = 5614
$ws.Range('G53').Value = 5
$ws.Range('H53').Value = 314

# Row 104: Maldivas
$ws.Range('D104').Value = 1863
$ws.Range('E104').Value = 412

# Row 118: Guinea-Bisau
$ws.Range('A118').Value = 'Guinea-Bisau'
$ws.Range('B118').Value = 1614
$ws.Range('C118').Value = 58
$ws.Range('D118').Value = 191
$ws.Range('E118').Value = 1401
$ws.Range('G118').Value = 3
$ws.Range('H118').Value = 22

# Row 119: Eslovenia
$ws.Range('A119').Value = 'Eslovenia'
$ws.Range('B119').Value = 1558
$ws.Range('C119').Value = 11
$ws.Range('D119').Value = 1376
$ws.Range('E119').Value = 73
$ws.Range('H119').Value = 109

# Row 120: Estado de Palestina
$ws.Range('A120').Value = 'Estado de Palestina'
$ws.Range('B120').Value = 1557
$ws.Range('C120').Value = 175
$ws.Range('D120').Value = 446
$ws.Range('E120').Value = 1108
$ws.Range('H120').Value = 3

# Row 121: Zambia
$ws.Range('A121').Value = 'Zambia'
$ws.Range('B121').Value = 1531
$ws.Range('C121').Value = 34
$ws.Range('D121').Value = 1233
$ws.Range('E121').Value = 277
$ws.Range('G121').Value = 3
$ws.Range('H121').Value = 21

# Row 122: Nueva Zelanda
$ws.Range('A122').Value = 'Nueva Zelanda'
$ws.Range('B122').Value = 1520
$ws.Range('C122').Value = 1
$ws.Range('D122').Value = 1484
$ws.Range('E122').Value = 14
$ws.Range('H122').Value = 22

# Row 132: Cabo Verde
$ws.Range('D132').Value = 568
$ws.Range('E132').Value = 449
$ws.Range('G132').Value = 2
$ws.Range('H132').Value = 10

# Row 133: Malaui
$ws.Range('A133').Value = 'Malaui'
$ws.Range('B133').Value = 1005
$ws.Range('C133').Value = 45
$ws.Range('D133').Value = 260
$ws.Range('E133').Value = 732
$ws.Range('G133').Value = 1
$ws.Range('H133').Value = 13

# Row 134: Republica de Chipre
$ws.Range('A134').Value = 'Republica de Chipre'
$ws.Range('B134').Value = 992
$ws.Range('D134').Value = 824
$ws.Range('E134').Value = 149
$ws.Range('H134').Value = 19

# Row 137: Uruguay
$ws.Range('B137').Value = 919
$ws.Range('C137').Value = 12
$ws.Range('E137').Value = 75

# Row 140: Principado de Andorra
$ws.Range('D140').Value = 799
$ws.Range('E140').Value = 4

# Row 143: Suazilandia
$ws.Range('A143').Value = 'Suazilandia'
$ws.Range('B143').Value = 728
$ws.Range('C143').Value = 22
$ws.Range('D143').Value = 367
$ws.Range('E143').Value = 353
$ws.Range('H143').Value = 8

# Row 144: Libia
$ws.Range('A144').Value = 'Libia'
$ws.Range('B144').Value = 713
$ws.Range('C144').Value = 15
$ws.Range('D144').Value = 142
$ws.Range('E144').Value = 553
$ws.Range('H144').Value = 18

# Row 145: Santo Tome y Principe
$ws.Range('A145').Value = 'Santo Tome y Principe'
$ws.Range('B145').Value = 712
$ws.Range('C145').Value = 1
$ws.Range('D145').Value = 219
$ws.Range('E145').Value = 480
$ws.Range('H145').Value = 13

# Row 146: Crucero
$ws.Range('A146').Value = 'Crucero'
$ws.Range('B146').Value = 712
$ws.Range('D146').Value = 651
$ws.Range('E146').Value = 48
$ws.Range('H146').Value = 13

# Row 151: Togo
$ws.Range('B151').Value = 591
$ws.Range('C151').Value = 3
$ws.Range('D151').Value = 395
$ws.Range('E151').Value = 182

# Row 153: Reunion
$ws.Range('B153').Value = 517
$ws.Range('C153').Value = 1
$ws.Range('E153').Value = 55

# Row 164: Martinica
$ws.Range('B164').Value = 242
$ws.Range('C164').Value = 6
$ws.Range('E164').Value = 130

# Row 168: Islas Caimanes
$ws.Range('D168').Value = 186
$ws.Range('E168').Value = 9

# Row 170: Guadalupe
$ws.Range('A170').Value = 'Guadalupe'
$ws.Range('B170').Value = 182
$ws.Range('C170').Value = 8
$ws.Range('D170').Value = 157
$ws.Range('E170').Value = 11
$ws.Range('H170').Value = 14

# Row 171: Gibraltar
$ws.Range('A171').Value = 'Gibraltar'
$ws.Range('B171').Value = 176
$ws.Range('D171').Value = 176
$ws.Range('E171').Value = 0
$ws.Range('H171').Value = 0

# Row 172: Eritrea
$ws.Range('A172').Value = 'Eritrea'
$ws.Range('B172').Value = 167
$ws.Range('C172').Value = 23
$ws.Range('D172').Value = 53
$ws.Range('E172').Value = 114
$ws.Range('H172').Value = 0

# Row 173: Bermudas
$ws.Range('A173').Value = 'Bermudas'
$ws.Range('B173').Value = 146
$ws.Range('D173').Value = 133
$ws.Range('E173').Value = 4
$ws.Range('H173').Value = 9

# Row 195: Lesoto
$ws.Range('A195').Value = 'Lesoto'
$ws.Range('C195').Value = 7
$ws.Range('D195').Value = 4
$ws.Range('E195').Value = 20

# Row 196: Timor Oriental
$ws.Range('A196').Value = 'Timor Oriental'
$ws.Range('B196').Value = 24
$ws.Range('D196').Value = 24
$ws.Range('E196').Value = 0
$ws.Range('H196').Value = 0

# Row 197: Belice
$ws.Range('A197').Value = 'Belice'
$ws.Range('D197').Value = 17
$ws.Range('E197').Value = 4
$ws.Range('H197').Value = 2

# Row 198: Curazao
$ws.Range('A198').Value = 'Curazao'
$ws.Range('D198').Value = 19
$ws.Range('E198').Value = 3
$ws.Range('H198').Value = 1

# Row 199: Granada
$ws.Range('A199').Value = 'Granada'
$ws.Range('B199').Value = 23
$ws.Range('D199').Value = 23

# Row 200: Nueva Caledonia
$ws.Range('A200').Value = 'Nueva Caledonia'
$ws.Range('B200').Value = 21
$ws.Range('D200').Value = 21

# Row 202: Laos
$ws.Range('A202').Value = 'Laos'
$ws.Range('B202').Value = 19
$ws.Range('D202').Value = 19

# Row 204: Dominica
$ws.Range('A204').Value = 'Dominica'
$ws.Range('B204').Value = 18
$ws.Range('D204').Value = 18
$ws.Range('E204').Value = 0

# Row 205: Islas Virgenes de los Estados Unidos
$ws.Range('A205').Value = 'Islas Virgenes de los Estados Unidos'
$ws.Range('D205').Value = 0
$ws.Range('E205').Value = 17

# Row 208: Groenlandia
$ws.Range('A208').Value = 'Groenlandia'

# Row 209: Islas Malvinas
$ws.Range('A209').Value = 'Islas Malvinas'

# Row 212: Seychelles
$ws.Range('A212').Value = 'Seychelles'
$ws.Range('D212').Value = 11
$ws.Range('H212').Value = 0

# Row 213: Montserrat
$ws.Range('A213').Value = 'Montserrat'
$ws.Range('D213').Value = 10
$ws.Range('H213').Value = 1
